$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$data = @(
    @('Bitcoin','https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc','66.483.58','  -1.26%  '),
    @('Ethereum','https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth','3.452.25','  -0.79%  '),
    @('TetherUSD','https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt','1.00','  +0.03%  '),
    @('BNB','https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb','580.02','  -2.33%  '),
    @('Solana','https://coinranking.com/coin/zNZHO_Sjf+solana-sol','175.74','  -1.60%  '),
    @('USDC','https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc','1.00','  +0.02%  '),
    @('XRP','https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp','0.597','  +0.61%  '),
    @('LidoStakedEther','https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth','3.451.34','  -0.83%  '),
    @('Dogecoin','https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge','0.134','  -2.77%  '),
    @('Toncoin','https://coinranking.com/coin/67YlI0K1b+toncoin-ton','6.85','  -3.22%  '),
    @('Cardano','https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada','0.418','  -3.26%  '),
    @('WrappedliquidstakedEther2.0','https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth','4.044.62','  -0.81%  '),
    @('Avalanche','https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax','30.76','  -3.92%  '),
    @('TRON','https://coinranking.com/coin/qUhEFk1I61atv+tron-trx','0.132','  -3.21%  '),
    @('WrappedBTC','https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc','66.466.92','  -1.28%  '),
    @('ShibaInu','https://coinranking.com/coin/xz24e0BjL+shibainu-shib','0.0000172','  -2.98%  '),
    @('WrappedEther','https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth','3.450.10','  -0.72%  '),
    @('Polkadot','https://coinranking.com/coin/25W7FG7om+polkadot-dot','5.99','  -4.20%  '),
    @('Chainlink','https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link','13.84','  -3.10%  '),
    @('BitcoinCash','https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch','376.11','  -3.30%  '),
    @('Uniswap','https://coinranking.com/coin/_H5FVG9iW+uniswap-uni','7.70','  -2.36%  '),
    @('Dai','https://coinranking.com/coin/MoTuySvg7+dai-dai','0.999','  -0.19%  '),
    @('LEO','https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo','5.71','  -0.07%  '),
    @('Polygon','https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic','0.527','  -1.47%  '),
    @('Litecoin','https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc','70.78','  -4.22%  '),
    @('PEPE','https://coinranking.com/coin/03WI8NQPF+pepe-pepe','0.0000117','  -3.17%  '),
    @('InternetComputer(DFINITY)','https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp','9.83','  -5.15%  '),
    @('Kaspa','https://coinranking.com/coin/V8GxkwWow+kaspa-kas','0.172','  -1.67%  '),
    @('Binance-PegBSC-USD','https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd','1.00','  -0.06%  '),
    @('NEARProtocol','https://coinranking.com/coin/DCrsaMv68+nearprotocol-near','5.85','  -5.15%  '),
    @('EthereumClassic','https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc','23.86','  +1.44%  '),
    @('PancakeSwap','https://coinranking.com/coin/ncYFcP709+pancakeswap-cake','1.98','  -3.71%  '),
    @('Fetch.AI','https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet','1.34','  -5.49%  '),
    @('USDe','https://coinranking.com/coin/exbfr2U-0+usde-usde','1.00','  -0.05%  '),
    @('Aptos','https://coinranking.com/coin/HGYj5JCv5+aptos-apt','7.02','  -4.48%  '),
    @('ImmutableX','https://coinranking.com/coin/Z96jIvLU7+immutablex-imx','1.51','  -4.85%  '),
    @('Monero','https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr','159.29','  -2.77%  '),
    @('Mantle','https://coinranking.com/coin/BoI4ux0nd+mantle-mnt','0.876','  +0.62%  '),
    @('EnergySwap','https://coinranking.com/coin/SbWqqTui-+energyswap-ens','27.21','  +3.65%  '),
    @('dogwifhat','https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif','2.63','  -3.52%  '),
    @('Stacks','https://coinranking.com/coin/mMPrMcB7+stacks-stx','1.78','  -5.04%  '),
    @('Filecoin','https://coinranking.com/coin/ymQub4fuB+filecoin-fil','4.44','  -4.05%  '),
    @('RenderToken','https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr','6.43','  -6.18%  '),
    @('Maker','https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr','2.690.92','  -5.03%  '),
    @('Hedera','https://coinranking.com/coin/jad286TjB+hedera-hbar','0.0695','  -3.15%  '),
    @('InjectiveProtocol','https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj','25.19','  -6.20%  '),
    @('OKB','https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb','40.24','  -3.08%  '),
    @('VeChain','https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet','0.0293','  -1.65%  '),
    @('Bittensor','https://coinranking.com/coin/pgv7xSFi6+bittensor-tao','321.21','  -4.48%  '),
    @('ONDO','https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo','1.01','  -3.85%  ')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    $ws.Cells.Item($rowNum, 2).Value = $rowVals[0]
    $ws.Cells.Item($rowNum, 3).Value = $rowVals[1]
    $ws.Cells.Item($rowNum, 4).Value = $rowVals[2]
    $ws.Cells.Item($rowNum, 5).Value = $rowVals[3]
}
